# Apply the change: insert a new "moment/2.22.2" row into the "3rd party"
# worksheet right before the existing "once/1.4.0" row (which is currently
# row 20), shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3rd party")

# Insert a new blank row at position 20; everything from row 20 downward
# (once, path-is-absolute, printj, ssf, wrappy, xlsx, xregexp) shifts down
# by one row, ending up on rows 21-27.
$ws.Rows.Item(20).Insert()

# Fill in the new row 20 with the moment/2.22.2 entry.
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = "moment/2.22.2"
$ws.Cells.Item(20, 3).Value = "Parse, validate, manipulate, and display dates"
$ws.Cells.Item(20, 4).Value = "MIT"
$ws.Cells.Item(20, 6).Value = "http://momentjs.com"
$ws.Cells.Item(20, 7).Value = "Iskren Ivov Chernev"

# Renumber column A for the shifted rows so the sequential index (1..26)
# is preserved rather than duplicated.
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(27, 1).Value = 26
